$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N (14th column), which pushes
# the existing "Late" / heading("Outstanding" label) / "Outstanding" columns
# one slot to the right (N->O, O->P, P->Q).
$ws.Columns("N:N").Insert()

# Give the freshly inserted column the same width as its left neighbour (M),
# matching the width Excel copies over automatically on a column insert.
$ws.Columns("N:N").ColumnWidth = $ws.Columns("M:M").ColumnWidth

# Update the sheet's selection / active cell to match the authored edit.
$ws.Range("R6").Select()

# Make "Repayment schedule" the active sheet/tab (it was "Acc_Periodic" before).
$ws.Activate()
